$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in missing "Value" column entries for R1, R2, R4
$ws.Range("D8").Value = "2k"
$ws.Range("D9").Value = "221R"
$ws.Range("D11").Value = "10K"

# Update the selection on Sheet1 to match the new active selection
$ws.Activate()
$ws.Range("A3:I17").Select()
